$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '47.367.04'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +5.73%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.650.93'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +12.48%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.34%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.35'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +8.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.79'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +14.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.617'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +11.51%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.614'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +24.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.98'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +18.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.13'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0862'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +11.94%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.54'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +24.16%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.046.83'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +12.08%  '

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.666.22'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +11.80%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.951'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +16.85%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '15.56'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +12.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '47.859.85'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +6.90%  '

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000105'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +13.70%  '

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.56'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +10.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.91'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +15.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.57'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +12.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '276.73'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +17.36%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.16'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +15.90%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '31.17'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +51.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.23'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +19.75%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.35%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.08'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.50%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.88'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +15.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '41.58'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +11.72%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.57%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.38'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +19.82%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.28%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.32'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +21.73%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0864'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +15.50%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.90'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +7.57%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '153.46'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.20%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +13.25%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +10.47%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.15'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +18.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.68'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +58.92%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.37'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +18.46%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.78'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +20.80%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0336'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +15.56%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.089.52'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +8.72%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.93'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +12.49%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.18%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '116.22'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +17.99%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +11.48%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.30'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +9.42%  '
